$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 6.712486666666666
$ws.Range("N2").Value = 20.13746
$ws.Range("O2").Value = 0.6330487633990675
$ws.Range("P2").Value = 0.6414503882251803
$ws.Range("Q2").Value = 196.0684173274244
$ws.Range("R2").Value = 1764.61575594682
$ws.Range("S2").Value = 0.01158161255404022
$ws.Range("T2").Value = 0.01180799862144473
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.290741083484562
$ws.Range("P3").Value = 0.2945997080427384
$ws.Range("Q3").Value = 90.04858296349468
$ws.Range("R3").Value = 810.437246671452
$ws.Range("S3").Value = 0.005319101429690936
$ws.Range("T3").Value = 0.005423074037061001
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 0.2495096666666667
$ws.Range("N4").Value = 0.748529
$ws.Range("O4").Value = 0.02353103905946135
$ws.Range("P4").Value = 0.02384333563656022
$ws.Range("Q4").Value = 7.288054022388112
$ws.Range("R4").Value = 65.592486201493
$ws.Range("S4").Value = 0.0004304998179245632
$ws.Range("T4").Value = 0.0004389148085265671
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 0.4166465
$ws.Range("N5").Value = 0.8332930000000001
$ws.Range("O5").Value = 0.03929356804674715
$ws.Range("P5").Value = 0.02654337331298611
$ws.Range("Q5").Value = 12.17003830274684
$ws.Range("R5").Value = 73.02022981648101
$ws.Range("S5").Value = 0.000718874922904416
$ws.Range("T5").Value = 0.0004886178592165817
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1419326666666667
$ws.Range("N6").Value = 0.425798
$ws.Range("O6").Value = 0.01338554601016197
$ws.Range("P6").Value = 0.01356319478253491
$ws.Range("Q6").Value = 4.145783031285112
$ws.Range("R6").Value = 37.312047281566
$ws.Range("S6").Value = 0.0002448882561298803
$ws.Range("T6").Value = 0.0002496750929369406
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 6.712486666666666
$ws.Range("N7").Value = 20.13746
$ws.Range("O7").Value = 0.6330487633990675
$ws.Range("P7").Value = 0.6414503882251803
$ws.Range("Q7").Value = 9788.691672541023
$ws.Range("R7").Value = 88098.22505286921
$ws.Range("S7").Value = 0.5782105854050416
$ws.Range("T7").Value = 0.5895128820369421
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.290741083484562
$ws.Range("P8").Value = 0.2945997080427384
$ws.Range("S8").Value = 0.2655554860896712
$ws.Range("T8").Value = 0.2707463057525654
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 0.2495096666666667
$ws.Range("N9").Value = 0.748529
$ws.Range("O9").Value = 0.02353103905946135
$ws.Range("P9").Value = 0.02384333563656022
$ws.Range("Q9").Value = 363.855202640028
$ws.Range("R9").Value = 3274.696823760253
$ws.Range("S9").Value = 0.02149265057671873
$ws.Range("T9").Value = 0.02191276794979259
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 0.4166465
$ws.Range("N10").Value = 0.8332930000000001
$ws.Range("O10").Value = 0.03929356804674715
$ws.Range("P10").Value = 0.02654337331298611
$ws.Range("Q10").Value = 607.5876686945668
$ws.Range("R10").Value = 3645.526012167401
$ws.Range("S10").Value = 0.03588974230195293
$ws.Range("T10").Value = 0.02439418665567602
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1419326666666667
$ws.Range("N11").Value = 0.425798
$ws.Range("O11").Value = 0.01338554601016197
$ws.Range("P11").Value = 0.01356319478253491
$ws.Range("Q11").Value = 206.9777090449651
$ws.Range("R11").Value = 1862.799381404686
$ws.Range("S11").Value = 0.01222601613333042
$ws.Range("T11").Value = 0.01246499837345752
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 6.712486666666666
$ws.Range("N12").Value = 20.13746
$ws.Range("O12").Value = 0.6330487633990675
$ws.Range("P12").Value = 0.6414503882251803
$ws.Range("Q12").Value = 389.2304487529422
$ws.Range("R12").Value = 3503.074038776479
$ws.Range("S12").Value = 0.0229915470994178
$ws.Range("T12").Value = 0.02344096343993999
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.290741083484562
$ws.Range("P13").Value = 0.2945997080427384
$ws.Range("Q13").Value = 178.7623465023253
$ws.Range("R13").Value = 1608.861118520928
$ws.Range("S13").Value = 0.01055935608937786
$ws.Range("T13").Value = 0.01076576008435214
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 0.2495096666666667
$ws.Range("N14").Value = 0.748529
$ws.Range("O14").Value = 0.02353103905946135
$ws.Range("P14").Value = 0.02384333563656022
$ws.Range("Q14").Value = 14.46807485028355
$ws.Range("R14").Value = 130.212673652552
$ws.Range("S14").Value = 0.0008546181970705396
$ws.Range("T14").Value = 0.0008713234401327102
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 0.4166465
$ws.Range("N15").Value = 0.8332930000000001
$ws.Range("O15").Value = 0.03929356804674715
$ws.Range("P15").Value = 0.02654337331298611
$ws.Range("Q15").Value = 24.15967617063066
$ws.Range("R15").Value = 144.958057023784
$ws.Range("S15").Value = 0.001427093729083645
$ws.Range("T15").Value = 0.0009699927770313595
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1419326666666667
$ws.Range("N16").Value = 0.425798
$ws.Range("O16").Value = 0.01338554601016197
$ws.Range("P16").Value = 0.01356319478253491
$ws.Range("Q16").Value = 8.230111772691556
$ws.Range("R16").Value = 74.07100595422399
$ws.Range("S16").Value = 0.0004861464540134606
$ws.Range("T16").Value = 0.0004956491707891448
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 6.712486666666666
$ws.Range("N17").Value = 20.13746
$ws.Range("O17").Value = 0.6330487633990675
$ws.Range("P17").Value = 0.6414503882251803
$ws.Range("Q17").Value = 197.8913899813666
$ws.Range("R17").Value = 1187.3483398882
$ws.Range("S17").Value = 0.01168929416468592
$ws.Range("T17").Value = 0.007945190058133583
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.290741083484562
$ws.Range("P18").Value = 0.2945997080427384
$ws.Range("Q18").Value = 90.88582185442002
$ws.Range("R18").Value = 545.31493112652
$ws.Range("S18").Value = 0.005368556495336101
$ws.Range("T18").Value = 0.003648997201399368
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 0.2495096666666667
$ws.Range("N19").Value = 0.748529
$ws.Range("O19").Value = 0.02353103905946135
$ws.Range("P19").Value = 0.02384333563656022
$ws.Range("Q19").Value = 7.355815691321666
$ws.Range("R19").Value = 44.13489414793
$ws.Range("S19").Value = 0.0004345024482629976
$ws.Range("T19").Value = 0.0002953304522528995
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 0.4166465
$ws.Range("N20").Value = 0.8332930000000001
$ws.Range("O20").Value = 0.03929356804674715
$ws.Range("P20").Value = 0.02654337331298611
$ws.Range("Q20").Value = 12.2831908814525
$ws.Range("R20").Value = 49.13276352581001
$ws.Range("S20").Value = 0.0007255587598217667
$ws.Range("T20").Value = 0.0003287738999413188
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.1419326666666667
$ws.Range("N21").Value = 0.425798
$ws.Range("O21").Value = 0.01338554601016197
$ws.Range("P21").Value = 0.01356319478253491
$ws.Range("Q21").Value = 4.184329010276667
$ws.Range("R21").Value = 25.10597406166
$ws.Range("S21").Value = 0.0002471651378443425
$ws.Range("T21").Value = 0.0001679976539431072
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 6.712486666666666
$ws.Range("N22").Value = 20.13746
$ws.Range("O22").Value = 0.6330487633990675
$ws.Range("P22").Value = 0.6414503882251803
$ws.Range("Q22").Value = 145.18087690778
$ws.Range("R22").Value = 1306.62789217002
$ws.Range("S22").Value = 0.008575724175881988
$ws.Range("T22").Value = 0.008743354068719875
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.290741083484562
$ws.Range("P23").Value = 0.2945997080427384
$ws.Range("Q23").Value = 66.67739974210801
$ws.Range("R23").Value = 600.096597678972
$ws.Range("S23").Value = 0.003938583380485842
$ws.Range("T23").Value = 0.004015570967360513
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 0.2495096666666667
$ws.Range("N24").Value = 0.748529
$ws.Range("O24").Value = 0.02353103905946135
$ws.Range("P24").Value = 0.02384333563656022
$ws.Range("Q24").Value = 5.396514585796999
$ws.Range("R24").Value = 48.568631272173
$ws.Range("S24").Value = 0.0003187680194845213
$ws.Range("T24").Value = 0.0003249989858554565
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 0.4166465
$ws.Range("N25").Value = 0.8332930000000001
$ws.Range("O25").Value = 0.03929356804674715
$ws.Range("P25").Value = 0.02654337331298611
$ws.Range("Q25").Value = 9.0114300756735
$ws.Range("R25").Value = 54.068580454041
$ws.Range("S25").Value = 0.0005322983329843907
$ws.Range("T25").Value = 0.0003618021211208262
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.1419326666666667
$ws.Range("N26").Value = 0.425798
$ws.Range("O26").Value = 0.01338554601016197
$ws.Range("P26").Value = 0.01356319478253491
$ws.Range("Q26").Value = 3.069787700414
$ws.Range("R26").Value = 27.628089303726
$ws.Range("S26").Value = 0.0001813300288438661
$ws.Range("T26").Value = 0.0001848744914081908
